# Updates Price (D) and Volume(1h) (E) columns for the cryptos list.
# D-column values that look numeric get forced to text (leading
# apostrophe) so Excel does not coerce them into numbers and drop
# formatting (e.g. "1.00" -> 1, "6.40" -> 6.4); the style is then
# reset to "Normal" so no stray quote-prefix formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.055.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.06%  '

$ws.Range("D3").Value = "'2.306.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.93%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'253.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.09%  '

$ws.Range("D6").Value = "'0.642"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.36%  '

$ws.Range("D7").Value = "'76.06"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.99%  '

$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  -3.13%  '

$ws.Range("D10").Value = "'39.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.55%  '

$ws.Range("D11").Value = "'0.0989"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.46%  '

$ws.Range("D12").Value = "'7.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.94%  '

$ws.Range("E13").Value = '  +2.29%  '

$ws.Range("D14").Value = "'2.652.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.72%  '

$ws.Range("D15").Value = "'15.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.37%  '

$ws.Range("D16").Value = "'0.884"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.40%  '

$ws.Range("D17").Value = "'2.294.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.49%  '

$ws.Range("D18").Value = "'43.048.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.05%  '

$ws.Range("E19").Value = '  +3.28%  '

$ws.Range("D21").Value = "'73.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.27%  '

$ws.Range("D22").Value = "'238.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("D23").Value = "'2.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +7.15%  '

$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("D25").Value = "'11.67"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("E27").Value = '  -1.51%  '

$ws.Range("E28").Value = '  -1.70%  '

$ws.Range("E29").Value = '  -1.06%  '

$ws.Range("D30").Value = "'167.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").Value = "'21.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.36%  '

$ws.Range("D32").Value = "'6.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.85%  '

$ws.Range("D33").Value = "'0.0845"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.68%  '

$ws.Range("E34").Value = '  -0.36%  '

$ws.Range("D35").Value = "'30.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.11%  '

$ws.Range("E36").Value = '  +1.77%  '

$ws.Range("D37").Value = "'4.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.14%  '

$ws.Range("E38").Value = '  +2.15%  '

$ws.Range("E39").Value = '  -2.45%  '

$ws.Range("D40").Value = "'13.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +12.06%  '

$ws.Range("D41").Value = "'2.36"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.05%  '

$ws.Range("D42").Value = "'5.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("E43").Value = '  +8.50%  '

$ws.Range("D44").Value = "'9.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.22%  '

$ws.Range("D45").Value = "'62.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.28%  '

$ws.Range("D46").Value = "'4.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.43%  '

$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("D48").Value = "'105.04"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.80%  '

$ws.Range("E49").Value = '  -0.55%  '

$ws.Range("E50").Value = '  -0.34%  '

$ws.Range("E51").Value = '  -0.68%  '
